$wb = $excel.ActiveWorkbook

$gamesWs = $wb.Worksheets.Item("Games")
$nextWs = $wb.Worksheets.Item("Next")

# The first upcoming game on the "Next" sheet (row 2) has now been played.
# Capture its data before removing it from "Next".
$playedDate = $nextWs.Cells.Item(2, 1).Value2
$playedOpp = $nextWs.Cells.Item(2, 2).Value2
$playedLoc = $nextWs.Cells.Item(2, 3).Value2

# Append the completed game to the bottom of the "Games" sheet.
$newRow = $gamesWs.Cells.Item($gamesWs.Rows.Count, 1).End(-4162).Row + 1
$newGameNum = $gamesWs.Cells.Item($newRow - 1, 1).Value2 + 1

$gamesWs.Cells.Item($newRow, 1).Value = $newGameNum
$gamesWs.Cells.Item($newRow, 2).Value = $playedDate
$gamesWs.Cells.Item($newRow, 2).NumberFormat = $gamesWs.Cells.Item($newRow - 1, 2).NumberFormat
$gamesWs.Cells.Item($newRow, 3).Value = -1
$gamesWs.Cells.Item($newRow, 4).Value = 94
$gamesWs.Cells.Item($newRow, 5).Value = 85.1
$gamesWs.Cells.Item($newRow, 6).Value = 0.476
$gamesWs.Cells.Item($newRow, 7).Value = 8.9
$gamesWs.Cells.Item($newRow, 8).Value = 28.3
$gamesWs.Cells.Item($newRow, 9).Value = 0.181
$gamesWs.Cells.Item($newRow, 10).Value = 110.5
$gamesWs.Cells.Item($newRow, 11).Value = $playedOpp
$gamesWs.Cells.Item($newRow, 12).Value = 98
$gamesWs.Cells.Item($newRow, 13).Value = 0.494
$gamesWs.Cells.Item($newRow, 14).Value = 6.9
$gamesWs.Cells.Item($newRow, 15).Value = 31.7
$gamesWs.Cells.Item($newRow, 16).Value = 0.165
$gamesWs.Cells.Item($newRow, 17).Value = 115.2
$gamesWs.Cells.Item($newRow, 18).Value = $playedLoc
$gamesWs.Cells.Item($newRow, 19).Value = 0

# Remove that now-played game from the "Next" sheet; remaining rows shift up.
$nextWs.Rows.Item(2).Delete()
